# Administracion de los mensajes de error mediante (Enum y tabla MENSAJES_SISTEMA)
# Adds a new "MENSAJES_SISTEMA" table box (Q13:T15) to the schema sheet,
# mirroring the existing table boxes already present (e.g. L13:O16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 13: table header "MENSAJES_SISTEMA" -------------------------------
# Copy format from the analogous header cell L13 (bold, yellow fill) onto Q13.
$ws.Range("L13").Copy() | Out-Null
$ws.Range("Q13").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Q13").Value = "MENSAJES_SISTEMA"

# Copy the header-row filler format (blue fill) from M13:O13 onto R13:T13.
$ws.Range("M13").Copy() | Out-Null
$ws.Range("R13:T13").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 14: ID / PK ---------------------------------------------------------
$ws.Range("L14").Copy() | Out-Null
$ws.Range("Q14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Q14").Value = "ID"

$ws.Range("M14").Copy() | Out-Null
$ws.Range("R14").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("R14").Value = "PK"

$ws.Range("N14").Copy() | Out-Null
$ws.Range("S14").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("O14").Copy() | Out-Null
$ws.Range("T14").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 15: DESCRIPCION / VACHAR(200) / NOT NULL / UNIQUE ------------------
$ws.Range("L15").Copy() | Out-Null
$ws.Range("Q15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("Q15").Value = "DESCRIPCION"

$ws.Range("M15").Copy() | Out-Null
$ws.Range("R15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("R15").Value = "VACHAR(200)"

$ws.Range("N15").Copy() | Out-Null
$ws.Range("S15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("S15").Value = "NOT NULL"

$ws.Range("O15").Copy() | Out-Null
$ws.Range("T15").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("T15").Value = "UNIQUE"

# --- Update the active selection to match the edited cell ------------------
$ws.Range("R15").Select() | Out-Null

$excel.CutCopyMode = $false
